# MAI_holdings.xlsx update:
#  - bump the "as of" date in the confidential disclaimer note from
#    2021-03-18 to 2021-03-19
#  - refresh the Weight / Percent Change figures for rows 2-7
#
# The worksheet ships protected (legacy password hash, no real secret is
# recoverable from it), so we briefly unprotect it to perform the writes and
# then re-protect it before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-19 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.4876501223183441
$ws.Range("E2").Value = 0.003977724741447641

$ws.Range("D3").Value = 0.3340203238596702
$ws.Range("E3").Value = -0.003476706069335589

$ws.Range("D4").Value = 0.0944933136466946
$ws.Range("E4").Value = -0.01739459255057685

$ws.Range("D5").Value = 0.05470188652277724
$ws.Range("E5").Value = 0.003849294296045747

$ws.Range("D6").Value = 0.02913435365251381
$ws.Range("E6").Value = 0.0134691195795007

$ws.Range("D7").Value = 0.9999999999999999
$ws.Range("E7").Value = -0.0002622474671745145

$ws.Protect()
